$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last refreshed" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 18 de Julio de 2020 a las 01:11"

# Refresh the per-country stats (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes) with the latest figures.
# Estados Unidos
$ws.Range("B4").Value = 3763088
$ws.Range("C4").Value = 68063
$ws.Range("D4").Value = 1710252
$ws.Range("E4").Value = 1910913
$ws.Range("G4").Value = 805
$ws.Range("H4").Value = 141923

# Brasil
$ws.Range("B5").Value = 2048697
$ws.Range("C5").Value = 33959
$ws.Range("E5").Value = 603990
$ws.Range("G5").Value = 1110
$ws.Range("H5").Value = 77932

# Peru
$ws.Range("B8").Value = 345537
$ws.Range("C8").Value = 3951
$ws.Range("D8").Value = 233982
$ws.Range("E8").Value = 98756
$ws.Range("G8").Value = 184
$ws.Range("H8").Value = 12799

# Colombia overtakes Francia in the ranking: row 21 now holds Colombia's
# figures and row 22 now holds Francia's (both rows keep their position,
# only the country label + stats move down/up one slot).
$ws.Range("A21").Value = "Colombia"
$ws.Range("B21").Value = 182140
$ws.Range("C21").Value = 8934
$ws.Range("D21").Value = 80637
$ws.Range("E21").Value = 95215
$ws.Range("G21").Value = 259
$ws.Range("H21").Value = 6288

$ws.Range("A22").Value = "Francia"
$ws.Range("B22").Value = 174674
$ws.Range("C22").Value = 836
$ws.Range("D22").Value = 79233
$ws.Range("E22").Value = 65289
$ws.Range("G22").Value = 14
$ws.Range("H22").Value = 30152

# Argentina
$ws.Range("B23").Value = 119301
$ws.Range("C23").Value = 4518
$ws.Range("E23").Value = 67343
$ws.Range("G23").Value = 66
$ws.Range("H23").Value = 2178

# Canada
$ws.Range("B24").Value = 109639
$ws.Range("C24").Value = 375
$ws.Range("D24").Value = 96674
$ws.Range("E24").Value = 4126
$ws.Range("G24").Value = 12
$ws.Range("H24").Value = 8839

# Nigeria overtakes Afganistan: row 50 now holds Nigeria's figures and row 51
# now holds Afganistan's.
$ws.Range("A50").Value = "Nigeria"
$ws.Range("B50").Value = 35454
$ws.Range("C50").Value = 600
$ws.Range("D50").Value = 14633
$ws.Range("E50").Value = 20049
$ws.Range("G50").Value = 3
$ws.Range("H50").Value = 772

$ws.Range("A51").Value = "Afganistan"
$ws.Range("B51").Value = 35229
$ws.Range("C51").Value = 159
$ws.Range("D51").Value = 23151
$ws.Range("E51").Value = 10931
$ws.Range("G51").Value = 34
$ws.Range("H51").Value = 1147

# Guatemala overtakes Suiza: row 53 now holds Guatemala's figures and row 54
# now holds Suiza's.
$ws.Range("A53").Value = "Guatemala"
$ws.Range("B53").Value = 33809
$ws.Range("C53").Value = 870
$ws.Range("D53").Value = 4989
$ws.Range("E53").Value = 27377
$ws.Range("G53").Value = 39
$ws.Range("H53").Value = 1443

$ws.Range("A54").Value = "Suiza"
$ws.Range("B54").Value = 33382
$ws.Range("C54").Value = 92
$ws.Range("D54").Value = 29900
$ws.Range("E54").Value = 1513
$ws.Range("H54").Value = 1969

# Ghana
$ws.Range("B57").Value = 26572
$ws.Range("C57").Value = 447
$ws.Range("D57").Value = 22915
$ws.Range("E57").Value = 3513
$ws.Range("G57").Value = 5
$ws.Range("H57").Value = 144

# Japon
$ws.Range("B59").Value = 23473
$ws.Range("C59").Value = 583
$ws.Range("D59").Value = 19096
$ws.Range("E59").Value = 3392

# Uzbekistan
$ws.Range("B67").Value = 15607
$ws.Range("C67").Value = 541
$ws.Range("D67").Value = 9003
$ws.Range("E67").Value = 6525
$ws.Range("G67").Value = 4
$ws.Range("H67").Value = 79

# Noruega
$ws.Range("B79").Value = 9025
$ws.Range("C79").Value = 10
$ws.Range("E79").Value = 632

# Paraguay
$ws.Range("B103").Value = 3457
$ws.Range("C103").Value = 115
$ws.Range("D103").Value = 1481
$ws.Range("E103").Value = 1948
$ws.Range("G103").Value = 1
$ws.Range("H103").Value = 28

# Surinam
$ws.Range("B147").Value = 943
$ws.Range("C147").Value = 39
$ws.Range("D147").Value = 600
$ws.Range("E147").Value = 324

# Birmania
$ws.Range("D164").Value = 271
$ws.Range("E164").Value = 62

# Bahamas
$ws.Range("B181").Value = 129
$ws.Range("C181").Value = 5
$ws.Range("E181").Value = 27
